$wb = $excel.ActiveWorkbook

# Update base year (2025) technology portfolio value in B2.
# Downstream sheets (2030, 2035, 2040, 2045, 2050) reference '2025'!B2
# via formulas, so they will recalculate automatically.
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 220000

$wb.Application.Calculate()
